# Change aspect (column K, "ASP") to 180 - abs(180 - ASP), i.e. fold the
# compass bearing (0-360) into the 0-180 range measuring deviation from
# north, symmetric around south (180).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column K (data starts at row 2, header in row 1)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 11).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 11)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = 180 - [Math]::Abs(180 - [double]$val)
    }
}
